$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.984.05"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "3.071.35"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'536.19"
$ws.Range("E5").Value = "  -4.69%  "
$ws.Range("D6").Value = "'133.72"
$ws.Range("E6").Value = "  -4.47%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.060.60"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("D11").Value = "'6.18"
$ws.Range("E11").Value = "  -8.90%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "'0.0000225"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "'34.39"
$ws.Range("E14").Value = "  -4.52%  "
$ws.Range("D15").Value = "3.567.61"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "62.933.41"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "3.067.32"
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("D19").Value = "'6.61"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "'484.16"
$ws.Range("E20").Value = "  -4.99%  "
$ws.Range("D21").Value = "'13.29"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").Value = "'0.696"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").Value = "'7.14"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").Value = "'79.27"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "'12.14"
$ws.Range("E25").Value = "  -3.96%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -4.42%  "
$ws.Range("D28").Value = "'8.12"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "'25.96"
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("D31").Value = "'1.87"
$ws.Range("E31").Value = "  -10.23%  "
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").Value = "'2.39"
$ws.Range("E33").Value = "  -7.42%  "
$ws.Range("D34").Value = "'56.74"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "'5.93"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'472.45"
$ws.Range("E37").Value = "  -15.02%  "
$ws.Range("D38").Value = "3.107.52"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").Value = "'0.0394"
$ws.Range("E39").Value = "  -6.76%  "
$ws.Range("D40").Value = "'0.0794"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").Value = "'0.115"
$ws.Range("E41").Value = "  -4.84%  "
$ws.Range("D42").Value = "'8.07"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").Value = "'2.63"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.252"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "0.0₃0547"
$ws.Range("E46").Value = "  +7.08%  "
$ws.Range("D47").Value = "'2.04"
$ws.Range("E47").Value = "  -5.36%  "
$ws.Range("D48").Value = "'120.73"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "'24.46"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "'0.108"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "'2.35"
$ws.Range("E51").Value = "  +4.58%  "
